$d = $word.ActiveDocument

# Locate the target paragraph: "En este paso del aprendizaje..." (last paragraph in the body).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt.StartsWith("En este paso del aprendizaje")) {
        $target = $i
    }
}

$p1 = $d.Paragraphs($target)

# Create two new empty paragraphs right after it; they inherit the
# paragraph-mark formatting (Arial/24/underline) of $p1 at the time of
# insertion, same as Word does when you press Enter at the end of a line.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($target + 1)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($target + 2)

# --- Paragraph 2: the "Bokeh" paragraph -----------------------------------
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertAfter("Una vez que hemos conseguido almacenar la información en la base de datos Neo4j, el próximo paso es representarla correctamente. Debido a la poca información y documentación en internet con respecto al grafo de Neo4j, me he visto obligado a utilizar la herramienta Bokeh para representar dicha información almacenada.")

# Re-fetch the paragraph range/start now that it holds text.
$p2 = $d.Paragraphs($target + 1)
$p2Start = $p2.Range.Start

# The final "información" (right before " almacenada.") is split into
# "informació" + "n" and both halves carry a single underline; the
# relocated _GoBack bookmark sits between the two halves.
$uStart = $p2Start + 294
$uMid   = $p2Start + 304
$uEnd   = $p2Start + 305

$rUnderline1 = $d.Range($uStart, $uMid)
$rUnderline1.Font.Underline = 1

$rUnderline2 = $d.Range($uMid, $uEnd)
$rUnderline2.Font.Underline = 1

$bmRange = $d.Range($uMid, $uMid)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Paragraph 3: the "Pycharm" paragraph ----------------------------------
$p3 = $d.Paragraphs($target + 2)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter("Debido al desconocimiento de dicha librería de Python he procedido a instalarla, que gracias al IDE Pycharm ha sido muy sencillo de realizar. Antes de ponernos con el problema real y cohesionarlo todo, he ido poco a poco implementando aplicaciones más sencillas para familiarizarme con la librería.")

# --- Strip the paragraph-mark underline from paragraphs 1 and 2 -----------
# (paragraph 3 is the one that keeps the underlined paragraph mark, matching
# the target layout)
$p1 = $d.Paragraphs($target)
$p1.Range.Font.Underline = 0

$p2 = $d.Paragraphs($target + 1)
$p2.Range.Font.Underline = 0

Write-Output "done"
